# Auto-generated: apply cryptos.xlsx data refresh per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'67.731.24"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.04%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.800.41"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.40%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.03%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'599.07"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.67%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'167.53"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.85%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.10%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.44%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.01%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'6.29"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.09%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.08%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.0000254"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -0.77%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'35.98"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.50%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'4.440.23"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.53%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'3.830.23"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.78%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = "'Chainlink"
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = "'18.48"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.45%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = "'WrappedBTC"
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = "'67.806.02"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.20%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'7.09"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.58%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.51%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'461.74"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.83%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'9.89"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -2.94%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.71%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.0000152"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.03%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'83.35"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.23%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.79%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  -1.45%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.02%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.62%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'3.950.68"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.50%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.55%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.42%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.80%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.43%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.03%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -1.22%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'3.742.36"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.19%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.05%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'3.42"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +3.16%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -0.20%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.54%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.60%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.06%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D44').Value = "'48.08"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.97%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  +1.05%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'EnergySwap"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'28.02"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +11.35%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = "'Arweave"
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'42.78"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -3.79%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.43%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'ONDO"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'1.37"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +9.16%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'Monero"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'147.70"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.04%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.93%  "
$ws.Range('E51').Style = 'Normal'
